$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = "Greater thann USD500 for ABCMGUG2"
$ws.Range("A18").Value = "Greater thann USD600 for ABCMGUG2"
$ws.Range("D18").Value = "ABCMGUG2"

$ws.Range("A19").Select()
